$wb = $excel.ActiveWorkbook

$ws = $wb.Sheets.Item("ALC")
$ws.Range("H17").Value = 2506454.2
$ws.Range("J17").Value = 2638369.8
$ws.Range("L17").Value = 7915109.399999999
$ws.Range("N17").Value = -7915445.399999999

$ws.Range("H40").Value = 2281.3333
$ws.Range("I40").Value = 1738
$ws.Range("J40").Value = 2490.3076
$ws.Range("K40").Value = 1738
$ws.Range("L40").Value = 2490.3076
$ws.Range("M40").Value = -1563
$ws.Range("N40").Value = -2840.3076

$ws.Range("H64").Value = 3517.1428
$ws.Range("I64").Value = 3081.818
$ws.Range("J64").Value = 4253.846
$ws.Range("K64").Value = 3081.818
$ws.Range("L64").Value = 4253.846
$ws.Range("M64").Value = -2833.818
$ws.Range("N64").Value = -4749.846

$ws.Range("H67").Value = 3517.1428
$ws.Range("I67").Value = 3081.818
$ws.Range("J67").Value = 4253.846
$ws.Range("K67").Value = 3081.818
$ws.Range("L67").Value = 4253.846
$ws.Range("M67").Value = -2223.818
$ws.Range("N67").Value = -5969.846

$ws.Range("H76").Value = 3412.0256
$ws.Range("I76").Value = 2716.25
$ws.Range("J76").Value = 5183.091
$ws.Range("K76").Value = 2716.25
$ws.Range("L76").Value = 5183.091
$ws.Range("M76").Value = -2401.25
$ws.Range("N76").Value = -5813.091

$ws.Range("H79").Value = 3412.0256
$ws.Range("I79").Value = 2716.25
$ws.Range("J79").Value = 5183.091
$ws.Range("K79").Value = 2716.25
$ws.Range("L79").Value = 5183.091
$ws.Range("M79").Value = -1624.25
$ws.Range("N79").Value = -7367.091

$ws.Range("H87").Value = 14464.088
$ws.Range("J87").Value = 14464.088
$ws.Range("L87").Value = 14464.088
$ws.Range("N87").Value = -16960.088

$ws.Range("H90").Value = 14464.088
$ws.Range("J90").Value = 14464.088
$ws.Range("L90").Value = 43392.264
$ws.Range("N90").Value = -55872.264

$ws.Range("H112").Value = 1144.7028
$ws.Range("J112").Value = 1151.5
$ws.Range("L112").Value = 3454.5
$ws.Range("N112").Value = -5670.5

$ws.Range("H125").Value = 2353504.2
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 2353504.2
$ws.Range("K125").Value = 0
$ws.Range("N125").Value = -21186457.8
$ws.Range("M125").ClearContents()

$ws.Range("H138").Value = 2726.47
$ws.Range("I138").Value = 1514.0426
$ws.Range("J138").Value = 3801.6416
$ws.Range("K138").Value = 4542.1278
$ws.Range("L138").Value = 11404.9248
$ws.Range("M138").Value = 597.8721999999998
$ws.Range("N138").Value = -21684.9248

$ws = $wb.Sheets.Item("ARM")
$ws.Range("H63").Value = 2427.5
$ws.Range("I63").Value = 1686.3334
$ws.Range("K63").Value = 1686.3334
$ws.Range("M63").Value = -1000.3334

$ws.Range("H66").Value = 2427.5
$ws.Range("I66").Value = 1686.3334
$ws.Range("K66").Value = 8431.666999999999
$ws.Range("M66").Value = -4999.666999999999

$ws.Range("H122").Value = 1887.6
$ws.Range("I122").Value = 1656
$ws.Range("K122").Value = 4968
$ws.Range("M122").Value = -2518

$ws = $wb.Sheets.Item("BSM")
$ws.Range("H105").Value = 2300.8333
$ws.Range("I105").Value = 2190.889
$ws.Range("J105").Value = 2630.6667
$ws.Range("K105").Value = 2190.889
$ws.Range("L105").Value = 2630.6667
$ws.Range("M105").Value = -443.8890000000001
$ws.Range("N105").Value = -6124.6667

$ws = $wb.Sheets.Item("CRP")
$ws.Range("H62").Value = 5252
$ws.Range("I62").Value = 4835
$ws.Range("J62").Value = 6503
$ws.Range("K62").Value = 4835
$ws.Range("L62").Value = 6503
$ws.Range("M62").Value = -4211
$ws.Range("N62").Value = -7751

$ws.Range("H65").Value = 5252
$ws.Range("I65").Value = 4835
$ws.Range("J65").Value = 6503
$ws.Range("K65").Value = 24175
$ws.Range("L65").Value = 32515
$ws.Range("M65").Value = -21055
$ws.Range("N65").Value = -38755

$ws.Range("H140").Value = 83605.14
$ws.Range("J140").Value = 83605.14
$ws.Range("L140").Value = 83605.14
$ws.Range("N140").Value = -93965.14

$ws = $wb.Sheets.Item("GSM")
$ws.Range("H70").Value = 5120.108
$ws.Range("I70").Value = 4811.1763
$ws.Range("J70").Value = 5382.7
$ws.Range("K70").Value = 4811.1763
$ws.Range("L70").Value = 5382.7
$ws.Range("M70").Value = -4541.1763
$ws.Range("N70").Value = -5922.7

$ws.Range("H73").Value = 5120.108
$ws.Range("I73").Value = 4811.1763
$ws.Range("J73").Value = 5382.7
$ws.Range("K73").Value = 4811.1763
$ws.Range("L73").Value = 5382.7
$ws.Range("M73").Value = -3875.1763
$ws.Range("N73").Value = -7254.7

$ws.Range("H80").Value = 2862.4
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 2862.4
$ws.Range("K80").Value = 0
$ws.Range("N80").Value = -4858.4
$ws.Range("M80").ClearContents()

$ws.Range("H83").Value = 2862.4
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 2862.4
$ws.Range("K83").Value = 0
$ws.Range("N83").Value = -24296
$ws.Range("M83").ClearContents()

$ws.Range("H102").Value = 3078695.5
$ws.Range("I102").Value = 4274876.5
$ws.Range("J102").Value = 2802
$ws.Range("K102").Value = 4274876.5
$ws.Range("L102").Value = 2802
$ws.Range("M102").Value = -4273254.5
$ws.Range("N102").Value = -6046

$ws.Range("H122").Value = 2185.9556
$ws.Range("I122").Value = 1860.3667
$ws.Range("K122").Value = 5581.1001
$ws.Range("M122").Value = -3131.1001

$ws.Range("H126").Value = 1947.7843
$ws.Range("I126").Value = 1757.9062
$ws.Range("J126").Value = 2267.5789
$ws.Range("K126").Value = 5273.7186
$ws.Range("L126").Value = 6802.736699999999
$ws.Range("M126").Value = -2803.7186
$ws.Range("N126").Value = -11742.7367

$ws.Range("H132").Value = 2105.2432
$ws.Range("I132").Value = 1729.5385
$ws.Range("J132").Value = 2993.2727
$ws.Range("K132").Value = 5188.6155
$ws.Range("L132").Value = 8979.8181
$ws.Range("M132").Value = -2658.6155
$ws.Range("N132").Value = -14039.8181

$ws = $wb.Sheets.Item("LTW")
$ws.Range("H40").Value = 19258.69
$ws.Range("I40").Value = 26269.85
$ws.Range("J40").Value = 3678.3333
$ws.Range("K40").Value = 26269.85
$ws.Range("L40").Value = 3678.3333
$ws.Range("M40").Value = -26133.85
$ws.Range("N40").Value = -3950.3333

$ws.Range("H68").Value = 2051.4092
$ws.Range("I68").Value = 1546.6
$ws.Range("J68").Value = 2472.0833
$ws.Range("K68").Value = 1546.6
$ws.Range("L68").Value = 2472.0833
$ws.Range("M68").Value = -797.5999999999999
$ws.Range("N68").Value = -3970.0833

$ws.Range("H71").Value = 2051.4092
$ws.Range("I71").Value = 1546.6
$ws.Range("J71").Value = 2472.0833
$ws.Range("K71").Value = 7733
$ws.Range("L71").Value = 12360.4165
$ws.Range("M71").Value = -3989
$ws.Range("N71").Value = -19848.4165

$ws.Range("H132").Value = 9861.633
$ws.Range("I132").Value = 8009.1113
$ws.Range("J132").Value = 14991.692
$ws.Range("K132").Value = 24027.3339
$ws.Range("L132").Value = 44975.076
$ws.Range("M132").Value = -21497.3339
$ws.Range("N132").Value = -50035.076
